$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name cells (column A) due to reordering of the countries list ---
$ws.Cells.Item(1, 1).Value2 = 'Datos actualizados a 2 de Abril de 2020 a las 01:50'
$ws.Cells.Item(46, 1).Value2 = 'Panama'
$ws.Cells.Item(47, 1).Value2 = 'Republica Dominicana'
$ws.Cells.Item(48, 1).Value2 = 'Islandia'
$ws.Cells.Item(49, 1).Value2 = 'Mexico'
$ws.Cells.Item(50, 1).Value2 = 'Argentina'
$ws.Cells.Item(51, 1).Value2 = 'Colombia'
$ws.Cells.Item(52, 1).Value2 = 'Serbia'
$ws.Cells.Item(133, 1).Value2 = 'Barbados'
$ws.Cells.Item(134, 1).Value2 = 'Uganda'
$ws.Cells.Item(135, 1).Value2 = 'Jamaica'
$ws.Cells.Item(136, 1).Value2 = 'Macao'
$ws.Cells.Item(137, 1).Value2 = 'Puerto Rico'
$ws.Cells.Item(138, 1).Value2 = 'Guatemala'
$ws.Cells.Item(139, 1).Value2 = 'Polinesia Francesa'
$ws.Cells.Item(140, 1).Value2 = 'Zambia'
$ws.Cells.Item(141, 1).Value2 = 'Togo'
$ws.Cells.Item(143, 1).Value2 = 'El Salvador'
$ws.Cells.Item(144, 1).Value2 = 'Guam'
$ws.Cells.Item(158, 1).Value2 = 'Birmania'
$ws.Cells.Item(159, 1).Value2 = 'Haiti'
$ws.Cells.Item(170, 1).Value2 = 'Surinam'
$ws.Cells.Item(171, 1).Value2 = 'Mozambique'
$ws.Cells.Item(172, 1).Value2 = 'Libia'
$ws.Cells.Item(173, 1).Value2 = 'Seychelles'
$ws.Cells.Item(174, 1).Value2 = 'Laos'
$ws.Cells.Item(184, 1).Value2 = 'Antigua y Barbuda'
$ws.Cells.Item(185, 1).Value2 = 'Republica del Chad'
$ws.Cells.Item(188, 1).Value2 = 'Islas Turcas y Caicos'
$ws.Cells.Item(189, 1).Value2 = 'Liberia'
$ws.Cells.Item(194, 1).Value2 = 'Nepal'
$ws.Cells.Item(196, 1).Value2 = 'Nicaragua'
$ws.Cells.Item(200, 1).Value2 = 'Republica de Africa Central'
$ws.Cells.Item(201, 1).Value2 = 'Islas Virgenes Britanicas'
$ws.Cells.Item(202, 1).Value2 = 'Belice'
$ws.Cells.Item(203, 1).Value2 = 'Burundi'
$ws.Cells.Item(204, 1).Value2 = 'Anguila'
$ws.Cells.Item(205, 1).Value2 = 'Bonaire, San Eustaquio y Saba'
$ws.Cells.Item(206, 1).Value2 = 'Sierra Leona'
$ws.Cells.Item(208, 1).Value2 = 'Timor Oriental'
$ws.Cells.Item(209, 1).Value2 = 'Papua Nueva Guinea'

# --- Update statistics cells (columns B-H) with refreshed case numbers ---
$ws.Cells.Item(4, 2).Value2 = 214639
$ws.Cells.Item(4, 3).Value2 = 26109
$ws.Cells.Item(4, 4).Value2 = 8878
$ws.Cells.Item(4, 5).Value2 = 200662
$ws.Cells.Item(4, 7).Value2 = 1046
$ws.Cells.Item(4, 8).Value2 = 5099
$ws.Cells.Item(9, 4).Value2 = 10935
$ws.Cells.Item(9, 5).Value2 = 42022
$ws.Cells.Item(18, 2).Value2 = 9731
$ws.Cells.Item(18, 3).Value2 = 1119
$ws.Cells.Item(18, 5).Value2 = 7881
$ws.Cells.Item(20, 2).Value2 = 6880
$ws.Cells.Item(20, 3).Value2 = 1163
$ws.Cells.Item(20, 5).Value2 = 6511
$ws.Cells.Item(20, 7).Value2 = 41
$ws.Cells.Item(20, 8).Value2 = 242
$ws.Cells.Item(22, 2).Value2 = 5048
$ws.Cells.Item(22, 3).Value2 = 285
$ws.Cells.Item(22, 5).Value2 = 4680
$ws.Cells.Item(22, 7).Value2 = 3
$ws.Cells.Item(22, 8).Value2 = 23
$ws.Cells.Item(43, 5).Value2 = 1312
$ws.Cells.Item(43, 7).Value2 = 2
$ws.Cells.Item(43, 8).Value2 = 51
$ws.Cells.Item(46, 2).Value2 = 1317
$ws.Cells.Item(46, 3).Value2 = 136
$ws.Cells.Item(46, 5).Value2 = 1276
$ws.Cells.Item(46, 6).Value2 = 50
$ws.Cells.Item(46, 7).Value2 = 2
$ws.Cells.Item(46, 8).Value2 = 32
$ws.Cells.Item(47, 2).Value2 = 1284
$ws.Cells.Item(47, 3).Value2 = 175
$ws.Cells.Item(47, 4).Value2 = 9
$ws.Cells.Item(47, 5).Value2 = 1218
$ws.Cells.Item(47, 6).Value2 = 0
$ws.Cells.Item(47, 7).Value2 = 6
$ws.Cells.Item(47, 8).Value2 = 57
$ws.Cells.Item(48, 2).Value2 = 1220
$ws.Cells.Item(48, 3).Value2 = 85
$ws.Cells.Item(48, 4).Value2 = 236
$ws.Cells.Item(48, 5).Value2 = 982
$ws.Cells.Item(48, 6).Value2 = 12
$ws.Cells.Item(48, 7).Value2 = 0
$ws.Cells.Item(48, 8).Value2 = 2
$ws.Cells.Item(49, 2).Value2 = 1215
$ws.Cells.Item(49, 3).Value2 = 121
$ws.Cells.Item(49, 4).Value2 = 35
$ws.Cells.Item(49, 5).Value2 = 1151
$ws.Cells.Item(49, 6).Value2 = 1
$ws.Cells.Item(49, 7).Value2 = 1
$ws.Cells.Item(49, 8).Value2 = 29
$ws.Cells.Item(50, 2).Value2 = 1133
$ws.Cells.Item(50, 3).Value2 = 79
$ws.Cells.Item(50, 4).Value2 = 248
$ws.Cells.Item(50, 5).Value2 = 853
$ws.Cells.Item(50, 6).Value2 = 0
$ws.Cells.Item(50, 7).Value2 = 5
$ws.Cells.Item(50, 8).Value2 = 32
$ws.Cells.Item(51, 2).Value2 = 1065
$ws.Cells.Item(51, 3).Value2 = 159
$ws.Cells.Item(51, 4).Value2 = 39
$ws.Cells.Item(51, 5).Value2 = 1009
$ws.Cells.Item(51, 6).Value2 = 47
$ws.Cells.Item(51, 7).Value2 = 1
$ws.Cells.Item(51, 8).Value2 = 17
$ws.Cells.Item(52, 2).Value2 = 1060
$ws.Cells.Item(52, 3).Value2 = 160
$ws.Cells.Item(52, 4).Value2 = 42
$ws.Cells.Item(52, 5).Value2 = 990
$ws.Cells.Item(52, 6).Value2 = 62
$ws.Cells.Item(52, 7).Value2 = 5
$ws.Cells.Item(52, 8).Value2 = 28
$ws.Cells.Item(83, 2).Value2 = 350
$ws.Cells.Item(83, 3).Value2 = 12
$ws.Cells.Item(83, 4).Value2 = 62
$ws.Cells.Item(83, 5).Value2 = 286
$ws.Cells.Item(83, 6).Value2 = 15
$ws.Cells.Item(108, 2).Value2 = 144
$ws.Cells.Item(108, 3).Value2 = 9
$ws.Cells.Item(108, 4).Value2 = 43
$ws.Cells.Item(108, 5).Value2 = 98
$ws.Cells.Item(133, 2).Value2 = 45
$ws.Cells.Item(133, 3).Value2 = 11
$ws.Cells.Item(133, 5).Value2 = 45
$ws.Cells.Item(134, 3).Value2 = 0
$ws.Cells.Item(134, 4).Value2 = 0
$ws.Cells.Item(134, 5).Value2 = 44
$ws.Cells.Item(134, 7).Value2 = 0
$ws.Cells.Item(134, 8).Value2 = 0
$ws.Cells.Item(135, 2).Value2 = 44
$ws.Cells.Item(135, 3).Value2 = 6
$ws.Cells.Item(135, 4).Value2 = 2
$ws.Cells.Item(135, 5).Value2 = 39
$ws.Cells.Item(135, 7).Value2 = 1
$ws.Cells.Item(135, 8).Value2 = 3
$ws.Cells.Item(136, 2).Value2 = 41
$ws.Cells.Item(136, 4).Value2 = 10
$ws.Cells.Item(136, 5).Value2 = 31
$ws.Cells.Item(136, 8).Value2 = 0
$ws.Cells.Item(137, 3).Value2 = 0
$ws.Cells.Item(137, 4).Value2 = 1
$ws.Cells.Item(137, 5).Value2 = 36
$ws.Cells.Item(137, 6).Value2 = 0
$ws.Cells.Item(137, 8).Value2 = 2
$ws.Cells.Item(138, 2).Value2 = 39
$ws.Cells.Item(138, 3).Value2 = 1
$ws.Cells.Item(138, 4).Value2 = 12
$ws.Cells.Item(138, 5).Value2 = 26
$ws.Cells.Item(138, 8).Value2 = 1
$ws.Cells.Item(139, 2).Value2 = 37
$ws.Cells.Item(139, 5).Value2 = 37
$ws.Cells.Item(139, 6).Value2 = 1
$ws.Cells.Item(140, 3).Value2 = 0
$ws.Cells.Item(140, 4).Value2 = 0
$ws.Cells.Item(140, 5).Value2 = 36
$ws.Cells.Item(140, 7).Value2 = 0
$ws.Cells.Item(140, 8).Value2 = 0
$ws.Cells.Item(141, 2).Value2 = 36
$ws.Cells.Item(141, 3).Value2 = 2
$ws.Cells.Item(141, 4).Value2 = 10
$ws.Cells.Item(141, 5).Value2 = 24
$ws.Cells.Item(141, 7).Value2 = 1
$ws.Cells.Item(141, 8).Value2 = 2
$ws.Cells.Item(143, 2).Value2 = 33
$ws.Cells.Item(143, 3).Value2 = 1
$ws.Cells.Item(143, 6).Value2 = 4
$ws.Cells.Item(143, 7).Value2 = 1
$ws.Cells.Item(143, 8).Value2 = 2
$ws.Cells.Item(144, 5).Value2 = 31
$ws.Cells.Item(144, 6).Value2 = 0
$ws.Cells.Item(144, 7).Value2 = 0
$ws.Cells.Item(144, 8).Value2 = 1
$ws.Cells.Item(157, 4).Value2 = 1
$ws.Cells.Item(157, 5).Value2 = 15
$ws.Cells.Item(158, 4).Value2 = 0
$ws.Cells.Item(158, 8).Value2 = 1
$ws.Cells.Item(159, 4).Value2 = 1
$ws.Cells.Item(159, 8).Value2 = 0
$ws.Cells.Item(170, 3).Value2 = 1
$ws.Cells.Item(171, 3).Value2 = 2
$ws.Cells.Item(172, 3).Value2 = 0
$ws.Cells.Item(173, 3).Value2 = 0
$ws.Cells.Item(188, 3).Value2 = 1
$ws.Cells.Item(189, 3).Value2 = 3
$ws.Cells.Item(194, 4).Value2 = 1
$ws.Cells.Item(194, 8).Value2 = 0
$ws.Cells.Item(196, 4).Value2 = 0
$ws.Cells.Item(196, 8).Value2 = 1
$ws.Cells.Item(203, 3).Value2 = 0
$ws.Cells.Item(204, 3).Value2 = 0
$ws.Cells.Item(205, 3).Value2 = 2
